$wb = $excel.ActiveWorkbook

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1581.6964
$ws.Range("I132").Value = 1648.8043
$ws.Range("J132").Value = 1273
$ws.Range("K132").Value = 4946.4129
$ws.Range("L132").Value = 3819
$ws.Range("M132").Value = -2416.4129
$ws.Range("N132").Value = -8879

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 346
$ws.Range("J135").Value = 550.5
$ws.Range("L135").Value = 4954.5
$ws.Range("N135").Value = -10024.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2665.9697
$ws.Range("I138").Value = 1239.2167
$ws.Range("K138").Value = 3717.6501
$ws.Range("M138").Value = 1422.3499

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4747.491
$ws.Range("I141").Value = 1161.8937
$ws.Range("J141").Value = 21599.8
$ws.Range("K141").Value = 3485.6811
$ws.Range("L141").Value = 64799.39999999999
$ws.Range("M141").Value = 1694.3189
$ws.Range("N141").Value = -75159.39999999999

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4850.94
$ws.Range("I32").Value = 4265.5137
$ws.Range("J32").Value = 9664.444
$ws.Range("K32").Value = 4265.5137
$ws.Range("L32").Value = 9664.444
$ws.Range("M32").Value = -3978.5137
$ws.Range("N32").Value = -10238.444

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7000.5557
$ws.Range("I61").Value = 8876.615
$ws.Range("J61").Value = 2122.8
$ws.Range("K61").Value = 8876.615
$ws.Range("L61").Value = 2122.8
$ws.Range("M61").Value = -8664.615
$ws.Range("N61").Value = -2546.8

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4372.3335
$ws.Range("I74").Value = 4819.0967
$ws.Range("J74").Value = 1602.4
$ws.Range("K74").Value = 4819.0967
$ws.Range("L74").Value = 1602.4
$ws.Range("M74").Value = -3945.0967
$ws.Range("N74").Value = -3350.4

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4372.3335
$ws.Range("I77").Value = 4819.0967
$ws.Range("J77").Value = 1602.4
$ws.Range("K77").Value = 24095.4835
$ws.Range("L77").Value = 8012
$ws.Range("M77").Value = -19727.4835
$ws.Range("N77").Value = -16748

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7000.5557
$ws.Range("I136").Value = 8876.615
$ws.Range("J136").Value = 2122.8
$ws.Range("K136").Value = 26629.845
$ws.Range("L136").Value = 6368.400000000001
$ws.Range("M136").Value = -24079.845
$ws.Range("N136").Value = -11468.4

# BSM row 7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# BSM row 113
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 2800
$ws.Range("I113").Value = 2800
$ws.Range("K113").Value = 2800
$ws.Range("M113").Value = -630

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2534.3262
$ws.Range("I31").Value = 1680.3846
$ws.Range("J31").Value = 3644.45
$ws.Range("K31").Value = 1680.3846
$ws.Range("L31").Value = 3644.45
$ws.Range("M31").Value = -1385.3846
$ws.Range("N31").Value = -4234.45

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2534.3262
$ws.Range("I34").Value = 1680.3846
$ws.Range("J34").Value = 3644.45
$ws.Range("K34").Value = 1680.3846
$ws.Range("L34").Value = 3644.45
$ws.Range("M34").Value = -1478.3846
$ws.Range("N34").Value = -4048.45

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1325.1852
$ws.Range("I58").Value = 977.70734
$ws.Range("J58").Value = 2421.077
$ws.Range("K58").Value = 977.70734
$ws.Range("L58").Value = 2421.077
$ws.Range("M58").Value = -774.70734
$ws.Range("N58").Value = -2827.077

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1505.5454
$ws.Range("I132").Value = 849.8919
$ws.Range("K132").Value = 2549.6757
$ws.Range("M132").Value = -19.67569999999978

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1215.4415
$ws.Range("I134").Value = 1187.9181
$ws.Range("J134").Value = 1320.375
$ws.Range("K134").Value = 3563.754300000001
$ws.Range("L134").Value = 3961.125
$ws.Range("M134").Value = -1028.754300000001
$ws.Range("N134").Value = -9031.125

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1325.1852
$ws.Range("I136").Value = 977.70734
$ws.Range("J136").Value = 2421.077
$ws.Range("K136").Value = 2933.12202
$ws.Range("L136").Value = 7263.231000000001
$ws.Range("M136").Value = -383.1220200000002
$ws.Range("N136").Value = -12363.231

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5222.353
$ws.Range("J131").Value = 5808.3335
$ws.Range("L131").Value = 17425.0005
$ws.Range("N131").Value = -27505.0005

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1826.3954
$ws.Range("I132").Value = 1526.7576
$ws.Range("J132").Value = 2815.2
$ws.Range("K132").Value = 4580.2728
$ws.Range("L132").Value = 8445.599999999999
$ws.Range("M132").Value = -2050.2728
$ws.Range("N132").Value = -13505.6

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 12177
$ws.Range("J136").Value = 12177
$ws.Range("L136").Value = 36531
$ws.Range("N136").Value = -41631

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 987.2727
$ws.Range("I22").Value = 533.3333
$ws.Range("J22").Value = 1157.5
$ws.Range("K22").Value = 533.3333
$ws.Range("L22").Value = 1157.5
$ws.Range("M22").Value = -238.3333
$ws.Range("N22").Value = -1747.5

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 987.2727
$ws.Range("I27").Value = 533.3333
$ws.Range("J27").Value = 1157.5
$ws.Range("K27").Value = 533.3333
$ws.Range("L27").Value = 1157.5
$ws.Range("M27").Value = -426.3333
$ws.Range("N27").Value = -1371.5

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8334858.5
$ws.Range("I136").Value = 1521.9706
$ws.Range("J136").Value = 55557100
$ws.Range("K136").Value = 4565.9118
$ws.Range("L136").Value = 166671300
$ws.Range("M136").Value = -2015.9118
$ws.Range("N136").Value = -166676400

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1366.7533
$ws.Range("I132").Value = 1285.8636
$ws.Range("J132").Value = 1852.091
$ws.Range("K132").Value = 3857.5908
$ws.Range("L132").Value = 5556.272999999999
$ws.Range("M132").Value = -1327.5908
$ws.Range("N132").Value = -10616.273

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3242.673
$ws.Range("I136").Value = 540.25714
$ws.Range("J136").Value = 8806.471
$ws.Range("K136").Value = 1620.77142
$ws.Range("L136").Value = 26419.413
$ws.Range("M136").Value = 929.22858
$ws.Range("N136").Value = -31519.413

Write-Output "Applied 162 cell updates"